$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update the raw input values in row 21 (J21 and K21).
# Dependent formulas (L21, M21, J24, K24, L24, M24, L25, M25, H27, I27)
# recalculate automatically.
$ws.Range("J21").Value = 11661
$ws.Range("K21").Value = 22

# Update the active cell / selection shown when the sheet was last saved.
$ws.Range("K22").Select()
